$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "276.17"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "23.27"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "6.469"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06300"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.658"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.690"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.384"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8366"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.01386"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1608"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08322"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03442"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03094"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09312"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.851"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001639"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.04788"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006296"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.005693"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.001088"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0001500"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.713"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04736"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007047"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1166"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003349"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01227"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006265"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.7960"
$ws.Range("E48").Value = "47CoinbaseStockTokenCOINWorstin24h"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002011"
